{"js": "// Update the date line and the changed multiplication problems in the\n// practice-sheet table. Every old text value in this document is unique,\n// so a literal (non-wildcard) search-and-replace for each pair is safe.\nconst replacements = [\n  [\"2024-10-11 Friday\", \"2024-10-12 Saturday\"],\n  [\"999\u00d73=\", \"549\u00d73=\"],\n  [\"447\u00d78=\", \"849\u00d79=\"],\n  [\"509\u00d78=\", \"285\u00d79=\"],\n  [\"673\u00d76=\", \"456\u00d78=\"],\n  [\"410\u00d77=\", \"561\u00d74=\"],\n  [\"821\u00d76=\", \"670\u00d75=\"],\n  [\"219\u00d75=\", \"800\u00d78=\"],\n  [\"391\u00d74=\", \"283\u00d76=\"],\n  [\"585\u00d72=\", \"945\u00d73=\"],\n  [\"441\u00d79=\", \"880\u00d76=\"],\n  [\"507\u00d75=\", \"551\u00d79=\"],\n  [\"563\u00d76=\", \"221\u00d76=\"],\n  [\"949\u00d76=\", \"152\u00d78=\"],\n  [\"948\u00d76=\", \"655\u00d77=\"],\n  [\"904\u00d72=\", \"106\u00d73=\"],\n  [\"763\u00d76=\", \"509\u00d77=\"],\n  [\"953\u00d77=\", \"551\u00d72=\"],\n  [\"866\u00d72=\", \"527\u00d79=\"],\n  [\"332\u00d79=\", \"361\u00d76=\"],\n  [\"976\u00d77=\", \"584\u00d72=\"],\n  [\"345\u00d74=\", \"396\u00d77=\"],\n  [\"898\u00d74=\", \"232\u00d72=\"],\n  [\"270\u00d75=\", \"884\u00d72=\"],\n  [\"174\u00d76=\", \"526\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the changed multiplication problems in the\n# practice-sheet table. Every old text value in this document is unique,\n# so a literal (non-wildcard) Find/Replace for each pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-11 Friday\", \"2024-10-12 Saturday\"),\n    @(\"999\u00d73=\", \"549\u00d73=\"),\n    @(\"447\u00d78=\", \"849\u00d79=\"),\n    @(\"509\u00d78=\", \"285\u00d79=\"),\n    @(\"673\u00d76=\", \"456\u00d78=\"),\n    @(\"410\u00d77=\", \"561\u00d74=\"),\n    @(\"821\u00d76=\", \"670\u00d75=\"),\n    @(\"219\u00d75=\", \"800\u00d78=\"),\n    @(\"391\u00d74=\", \"283\u00d76=\"),\n    @(\"585\u00d72=\", \"945\u00d73=\"),\n    @(\"441\u00d79=\", \"880\u00d76=\"),\n    @(\"507\u00d75=\", \"551\u00d79=\"),\n    @(\"563\u00d76=\", \"221\u00d76=\"),\n    @(\"949\u00d76=\", \"152\u00d78=\"),\n    @(\"948\u00d76=\", \"655\u00d77=\"),\n    @(\"904\u00d72=\", \"106\u00d73=\"),\n    @(\"763\u00d76=\", \"509\u00d77=\"),\n    @(\"953\u00d77=\", \"551\u00d72=\"),\n    @(\"866\u00d72=\", \"527\u00d79=\"),\n    @(\"332\u00d79=\", \"361\u00d76=\"),\n    @(\"976\u00d77=\", \"584\u00d72=\"),\n    @(\"345\u00d74=\", \"396\u00d77=\"),\n    @(\"898\u00d74=\", \"232\u00d72=\"),\n    @(\"270\u00d75=\", \"884\u00d72=\"),\n    @(\"174\u00d76=\", \"526\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
